$wb = $excel.ActiveWorkbook

# The change applies to both the "展览" and "全部类型" sheets, which carry
# the same data table (rows match in both).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F3: 想去人数 13 -> 14
    $ws.Range("F3").Value = 14

    # F4: 想去人数 966 -> 965
    $ws.Range("F4").Value = 965

    # F6: 想去人数 435 -> 436
    $ws.Range("F6").Value = 436
}
